$wb = $excel.ActiveWorkbook

# --- Commodities sheet: re-populate rows 4-46 with Csets (B) + CommName (C) ---
$wsComm = $wb.Worksheets.Item("Commodities")
$wsComm.Range("B4").Value = "NRG"
$wsComm.Range("C4").Value = "sec_biogas"
$wsComm.Range("B5").Value = "ENV"
$wsComm.Range("C5").Value = "emi_CO2_f_ind"
$wsComm.Range("B6").Value = "ENV"
$wsComm.Range("C6").Value = "emi_CO2_f_x2x_neg_reusable"
$wsComm.Range("B7").Value = "NRG"
$wsComm.Range("C7").Value = "pri_uran"
$wsComm.Range("B8").Value = "ENV"
$wsComm.Range("C8").Value = "emi_N2O_f_ind"
$wsComm.Range("B9").Value = "MAT"
$wsComm.Range("C9").Value = "iip_steel_oxygen"
$wsComm.Range("B10").Value = "MAT"
$wsComm.Range("C10").Value = "iip_steel_iron_pellets"
$wsComm.Range("B11").Value = "DEM"
$wsComm.Range("C11").Value = "exo_steel"
$wsComm.Range("B12").Value = "NRG"
$wsComm.Range("C12").Value = "pri_deuterium"
$wsComm.Range("B13").Value = "NRG"
$wsComm.Range("C13").Value = "pri_coal"
$wsComm.Range("B14").Value = "NRG"
$wsComm.Range("C14").Value = "iip_steel_blafu_slag"
$wsComm.Range("B15").Value = "NRG"
$wsComm.Range("C15").Value = "pri_crude_oil"
$wsComm.Range("B16").Value = "NRG"
$wsComm.Range("C16").Value = "pri_hydro_energy"
$wsComm.Range("B17").Value = "NRG"
$wsComm.Range("C17").Value = "sec_methane"
$wsComm.Range("B18").Value = "NRG"
$wsComm.Range("C18").Value = "sec_heating_oil"
$wsComm.Range("B19").Value = "NRG"
$wsComm.Range("C19").Value = "sec_hydrogen"
$wsComm.Range("B20").Value = "NRG"
$wsComm.Range("C20").Value = "sec_heat_low"
$wsComm.Range("B21").Value = "NRG"
$wsComm.Range("C21").Value = "iip_heat_proc"
$wsComm.Range("B22").Value = "NRG"
$wsComm.Range("C22").Value = "sec_elec_ind"
$wsComm.Range("B23").Value = "NRG"
$wsComm.Range("C23").Value = "pri_biomass"
$wsComm.Range("B24").Value = "NRG"
$wsComm.Range("C24").Value = "pri_waste"
$wsComm.Range("B25").Value = "NRG"
$wsComm.Range("C25").Value = "CO2_f_pow"
$wsComm.Range("B26").Value = "ENV"
$wsComm.Range("C26").Value = "emi_CO2_f_x2x_neg_stored]"
$wsComm.Range("B27").Value = "NRG"
$wsComm.Range("C27").Value = "iip_coke"
$wsComm.Range("B28").Value = "NRG"
$wsComm.Range("C28").Value = "pri_solar_radiation"
$wsComm.Range("B29").Value = "NRG"
$wsComm.Range("C29").Value = "sec_heavy_fuel_oil"
$wsComm.Range("B30").Value = "ENV"
$wsComm.Range("C30").Value = "emi_CH4_f_ind"
$wsComm.Range("B31").Value = "MAT"
$wsComm.Range("C31").Value = "iip_steel_raw_iron"
$wsComm.Range("B32").Value = "MAT"
$wsComm.Range("C32").Value = "iip_steel_crudesteel"
$wsComm.Range("B33").Value = "NRG"
$wsComm.Range("C33").Value = "pri_geoth_heat"
$wsComm.Range("B34").Value = "NRG"
$wsComm.Range("C34").Value = "pri_wind_energy_on"
$wsComm.Range("B35").Value = "MAT"
$wsComm.Range("C35").Value = "iip_steel_scrap"
$wsComm.Range("B36").Value = "NRG"
$wsComm.Range("C36").Value = "sec_elec"
$wsComm.Range("B37").Value = "NRG"
$wsComm.Range("C37").Value = "sec_natural_gas_syn"
$wsComm.Range("B38").Value = "ENV"
$wsComm.Range("C38").Value = "[emi_CO2_f_x2x_neg_reusable"
$wsComm.Range("B39").Value = "MAT"
$wsComm.Range("C39").Value = "iip_steel_iron_ore"
$wsComm.Range("B40").Value = "NRG"
$wsComm.Range("C40").Value = "pri_natural_gas"
$wsComm.Range("B41").Value = "NRG"
$wsComm.Range("C41").Value = "pri_wind_energy_off"
$wsComm.Range("B42").Value = "MAT"
$wsComm.Range("C42").Value = "iip_steel_sinter"
$wsComm.Range("B43").Value = "MAT"
$wsComm.Range("C43").Value = "iip_steel_sponge_iron"
$wsComm.Range("B44").Value = "NRG"
$wsComm.Range("C44").Value = "sec_heat_high"
$wsComm.Range("B45").Value = "ENV"
$wsComm.Range("C45").Value = "emi_CO2_f_x2x"
$wsComm.Range("B46").Value = "NRG"
$wsComm.Range("C46").Value = "sec_H2"

# --- Processes sheet: add Sets (B) for rows 4-77 ---
$wsProc = $wb.Worksheets.Item("Processes")
$wsProc.Range("B4").Value = "PRE"
$wsProc.Range("B5").Value = "PRE"
$wsProc.Range("B6").Value = "PRE"
$wsProc.Range("B7").Value = "PRE"
$wsProc.Range("B8").Value = "PRE"
$wsProc.Range("B9").Value = "PRE"
$wsProc.Range("B10").Value = "PRE"
$wsProc.Range("B11").Value = "PRE"
$wsProc.Range("B12").Value = "DEM"
$wsProc.Range("B13").Value = "DEM"
$wsProc.Range("B14").Value = "PRE"
$wsProc.Range("B15").Value = "PRE"
$wsProc.Range("B16").Value = "PRE"
$wsProc.Range("B17").Value = "PRE"
$wsProc.Range("B18").Value = "PRE"
$wsProc.Range("B19").Value = "PRE"
$wsProc.Range("B20").Value = "PRE"
$wsProc.Range("B21").Value = "PRE"
$wsProc.Range("B22").Value = "PRE"
$wsProc.Range("B23").Value = "PRE"
$wsProc.Range("B24").Value = "PRE"
$wsProc.Range("B25").Value = "PRE"
$wsProc.Range("B26").Value = "PRE"
$wsProc.Range("B27").Value = "PRE"
$wsProc.Range("B28").Value = "PRE"
$wsProc.Range("B29").Value = "PRE"
$wsProc.Range("B30").Value = "PRE"
$wsProc.Range("B31").Value = "PRE"
$wsProc.Range("B32").Value = "PRE"
$wsProc.Range("B33").Value = "PRE"
$wsProc.Range("B34").Value = "PRE"
$wsProc.Range("B35").Value = "CHP"
$wsProc.Range("B36").Value = "PRE"
$wsProc.Range("B37").Value = "CHP"
$wsProc.Range("B38").Value = "PRE"
$wsProc.Range("B39").Value = "PRE"
$wsProc.Range("B40").Value = "PRE"
$wsProc.Range("B41").Value = "PRE"
$wsProc.Range("B42").Value = "PRE"
$wsProc.Range("B43").Value = "PRE"
$wsProc.Range("B44").Value = "PRE"
$wsProc.Range("B45").Value = "PRE"
$wsProc.Range("B46").Value = "PRE"
$wsProc.Range("B47").Value = "PRE"
$wsProc.Range("B48").Value = "PRE"
$wsProc.Range("B49").Value = "PRE"
$wsProc.Range("B50").Value = "PRE"
$wsProc.Range("B51").Value = "PRE"
$wsProc.Range("B52").Value = "PRE"
$wsProc.Range("B53").Value = "PRE"
$wsProc.Range("B54").Value = "PRE"
$wsProc.Range("B55").Value = "PRE"
$wsProc.Range("B56").Value = "PRE"
$wsProc.Range("B57").Value = "PRE"
$wsProc.Range("B58").Value = "PRE"
$wsProc.Range("B59").Value = "PRE"
$wsProc.Range("B60").Value = "PRE"
$wsProc.Range("B61").Value = "PRE"
$wsProc.Range("B62").Value = "PRE"
$wsProc.Range("B63").Value = "PRE"
$wsProc.Range("B64").Value = "PRE"
$wsProc.Range("B65").Value = "PRE"
$wsProc.Range("B66").Value = "PRE"
$wsProc.Range("B67").Value = "PRE"
$wsProc.Range("B68").Value = "PRE"
$wsProc.Range("B69").Value = "PRE"
$wsProc.Range("B70").Value = "PRE"
$wsProc.Range("B71").Value = "PRE"
$wsProc.Range("B72").Value = "PRE"
$wsProc.Range("B73").Value = "PRE"
$wsProc.Range("B74").Value = "PRE"
$wsProc.Range("B75").Value = "PRE"
$wsProc.Range("B76").Value = "PRE"
$wsProc.Range("B77").Value = "PRE"
